$wb = $excel.ActiveWorkbook

# Rabbit sheet was previously the active tab with A6 selected; the user
# clicked on A5 there before moving on to edit the other sheets.
$rabbit = $wb.Worksheets.Item("rabbit")
$rabbit.Activate()
$rabbit.Range("A5").Select()

# Add a "healthy substitute" row for beef, reusing the plant-based
# ground beef substitute value, then leave the selection on B6.
$beef = $wb.Worksheets.Item("beef")
$beef.Activate()
$beef.Range("A6").Value = "healthy substitute"
$beef.Range("B6").Value = "plant-based ground beef"
$beef.Range("B6").Select()

# Add a "healthy substitute" row for pork, reusing the plant-based pork
# substitute value. Pork ends up the active sheet/selection.
$pork = $wb.Worksheets.Item("pork")
$pork.Activate()
$pork.Range("A6").Value = "healthy substitute"
$pork.Range("B6").Value = "plant-based pork"
$pork.Range("B13").Select()
